$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Wed Oct  2 14:11:33 UTC 2024 with GitHub Actions

# Row 2: Bitcoin
$ws.Range('D2').Value = '61.506.42'
$ws.Range('E2').Value = '  -2.09%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.452.80'
$ws.Range('E3').Value = '  -5.06%  '

# Row 4: TetherUSD
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  -0.09%  '

# Row 5: BNB
$ws.Range('D5').Value = "'546.50"
$ws.Range('E5').Value = '  -3.79%  '

# Row 6: Solana
$ws.Range('D6').Value = "'146.17"
$ws.Range('E6').Value = '  -4.67%  '

# Row 7: USDC
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  -0.10%  '

# Row 8: XRP
$ws.Range('E8').Value = '  -5.54%  '

# Row 9: LidoStakedEther
$ws.Range('D9').Value = '2.451.29'
$ws.Range('E9').Value = '  -5.08%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  -7.27%  '

# Row 11: TRON
$ws.Range('E11').Value = '  -1.20%  '

# Row 12: Toncoin
$ws.Range('E12').Value = '  -5.07%  '

# Row 13: Cardano
$ws.Range('D13').Value = "'0.351"
$ws.Range('E13').Value = '  -6.53%  '

# Row 14: Avalanche
$ws.Range('D14').Value = "'26.04"
$ws.Range('E14').Value = '  -6.66%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '2.894.52'
$ws.Range('E15').Value = '  -5.24%  '

# Row 16: ShibaInu
$ws.Range('D16').Value = "'0.0000168"
$ws.Range('E16').Value = '  -6.33%  '

# Row 17: WrappedBTC
$ws.Range('D17').Value = '61.331.12'
$ws.Range('E17').Value = '  -2.27%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '2.455.13'
$ws.Range('E18').Value = '  -4.95%  '

# Row 19: Chainlink
$ws.Range('D19').Value = "'10.97"
$ws.Range('E19').Value = '  -7.76%  '

# Row 20: Uniswap
$ws.Range('D20').Value = "'6.95"
$ws.Range('E20').Value = '  -6.54%  '

# Row 21: Polkadot
$ws.Range('D21').Value = "'4.16"
$ws.Range('E21').Value = '  -6.26%  '

# Row 22: BitcoinCash
$ws.Range('D22').Value = "'319.09"
$ws.Range('E22').Value = '  -5.17%  '

# Row 23: Dai
$ws.Range('E23').Value = '  +0.12%  '

# Row 24: SuiNetwork
$ws.Range('D24').Value = "'1.89"
$ws.Range('E24').Value = '  +1.45%  '

# Row 25: Litecoin
$ws.Range('D25').Value = "'63.74"
$ws.Range('E25').Value = '  -5.19%  '

# Row 26: PEPE
$ws.Range('D26').Value = '0.0₃0973'
$ws.Range('E26').Value = '  -11.00%  '

# Row 27: WrappedeETH
$ws.Range('D27').Value = '2.564.89'
$ws.Range('E27').Value = '  -5.71%  '

# Row 28: Binance-PegBSC-USD
$ws.Range('D28').Value = "'0.998"
$ws.Range('E28').Value = '  -0.08%  '

# Row 29: Fetch.AI -> Bittensor (swap)
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = "'535.25"
$ws.Range('E29').Value = '  -5.59%  '

# Row 30: Aptos
$ws.Range('D30').Value = "'7.85"
$ws.Range('E30').Value = '  -2.62%  '

# Row 31: Bittensor -> Fetch.AI (swap)
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = "'1.46"
$ws.Range('E31').Value = '  -9.53%  '

# Row 32: InternetComputer(DFINITY)
$ws.Range('D32').Value = "'8.25"
$ws.Range('E32').Value = '  -8.62%  '

# Row 33: Kaspa
$ws.Range('D33').Value = "'0.148"
$ws.Range('E33').Value = '  -6.39%  '

# Row 34: PancakeSwap
$ws.Range('D34').Value = "'1.88"
$ws.Range('E34').Value = '  -5.89%  '

# Row 35: ImmutableX
$ws.Range('E35').Value = '  -6.96%  '

# Row 36: RenderToken
$ws.Range('D36').Value = "'5.73"
$ws.Range('E36').Value = '  -10.71%  '

# Row 37: FirstDigitalUSD
$ws.Range('D37').Value = "'1.00"
$ws.Range('E37').Value = '  -0.16%  '

# Row 38: NEARProtocol
$ws.Range('D38').Value = "'4.81"
$ws.Range('E38').Value = '  -8.69%  '

# Row 39: PolygonEcosystemToken
$ws.Range('D39').Value = "'0.379"
$ws.Range('E39').Value = '  -4.30%  '

# Row 40: EthereumClassic
$ws.Range('D40').Value = "'18.25"
$ws.Range('E40').Value = '  -5.99%  '

# Row 41: Stacks
$ws.Range('E41').Value = '  -5.00%  '

# Row 42: Monero
$ws.Range('D42').Value = "'140.48"
$ws.Range('E42').Value = '  -8.49%  '

# Row 43: USDe
$ws.Range('E43').Value = '  +0.02%  '

# Row 44: OKB
$ws.Range('D44').Value = "'40.29"
$ws.Range('E44').Value = '  -3.26%  '

# Row 45: dogwifhat
$ws.Range('D45').Value = "'2.29"
$ws.Range('E45').Value = '  -7.36%  '

# Row 46: Aave
$ws.Range('D46').Value = "'140.89"
$ws.Range('E46').Value = '  -10.23%  '

# Row 47: Filecoin
$ws.Range('D47').Value = "'3.61"
$ws.Range('E47').Value = '  -6.17%  '

# Row 48: InjectiveProtocol
$ws.Range('D48').Value = "'21.33"
$ws.Range('E48').Value = '  -9.42%  '

# Row 49: Hedera
$ws.Range('D49').Value = "'0.0534"
$ws.Range('E49').Value = '  -7.39%  '

# Row 50: Mantle
$ws.Range('D50').Value = "'0.588"
$ws.Range('E50').Value = '  -5.43%  '

# Row 51: Stellar
$ws.Range('D51').Value = "'0.0931"
$ws.Range('E51').Value = '  -5.58%  '
